$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($hdrType in 1,2,3) {
        $hdr = $sec.Headers.Item($hdrType)
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("Dr. med. Thiên-Trí Lâm", $true, $false, $false, $false, $false,
                                     $true, 1, $false, "PD Dr. med. Thiên-Trí Lâm", 2)
        }
    }
}
